# Applies the "Updated cryptos list" data refresh described in the commit diff.
# Only cell VALUES are updated (coin name/link/price/volume-%); all cells are plain
# text in the source data, so values that look numeric (pure "123.45" style price
# strings) are written with a leading apostrophe to force Excel to keep them as text
# instead of silently coercing them to numbers (which would also normalize formatting,
# e.g. turn "93.10" into 93.1). The apostrophe itself is a text-qualifier and is not
# part of the stored cell value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.038.53"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").Value = "2.268.52"
$ws.Range("E3").Value = "  +0.58%  "

# Row 5
$ws.Range("D5").Value = "'305.71"
$ws.Range("E5").Value = "  +1.41%  "

# Row 6
$ws.Range("D6").Value = "'93.10"
$ws.Range("E6").Value = "  +1.33%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +1.58%  "

# Row 10
$ws.Range("D10").Value = "'32.79"
$ws.Range("E10").Value = "  +1.76%  "

# Row 11
$ws.Range("D11").Value = "'0.0799"
$ws.Range("E11").Value = "  +0.28%  "

# Row 12
$ws.Range("E12").Value = "  -1.73%  "

# Row 13
$ws.Range("E13").Value = "  +0.82%  "

# Row 14
$ws.Range("D14").Value = "2.621.27"
$ws.Range("E14").Value = "  +0.64%  "

# Row 15
$ws.Range("D15").Value = "'14.34"
$ws.Range("E15").Value = "  +1.89%  "

# Row 16
$ws.Range("D16").Value = "2.273.76"
$ws.Range("E16").Value = "  +0.54%  "

# Row 17
$ws.Range("E17").Value = "  +3.92%  "

# Row 18
$ws.Range("D18").Value = "41.909.04"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19
$ws.Range("D19").Value = "'12.77"
$ws.Range("E19").Value = "  +5.73%  "

# Row 20
$ws.Range("E20").Value = "  +2.01%  "

# Row 21
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").Value = "'68.18"
$ws.Range("E22").Value = "  +1.79%  "

# Row 23
$ws.Range("D23").Value = "'244.13"
$ws.Range("E23").Value = "  +1.47%  "

# Row 24
$ws.Range("E24").Value = "  +1.55%  "

# Row 25
$ws.Range("E25").Value = "  +2.65%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").Value = "'23.97"
$ws.Range("E27").Value = "  +0.46%  "

# Row 28
$ws.Range("E28").Value = "  +0.44%  "

# Row 29
$ws.Range("E29").Value = "  -9.46%  "

# Row 30
$ws.Range("D30").Value = "'34.92"
$ws.Range("E30").Value = "  +3.65%  "

# Row 31
$ws.Range("D31").Value = "'159.93"
$ws.Range("E31").Value = "  +0.69%  "

# Row 32
$ws.Range("D32").Value = "'5.33"
$ws.Range("E32").Value = "  +3.86%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("D34").Value = "'0.0743"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("E36").Value = "  +4.40%  "

# Row 37
$ws.Range("E37").Value = "  -1.03%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.105"
$ws.Range("E38").Value = "  +1.49%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.117"
$ws.Range("E39").Value = "  +1.14%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.79"
$ws.Range("E40").Value = "  +0.60%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'3.99"
$ws.Range("E41").Value = "  +1.88%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'19.78"
$ws.Range("E42").Value = "  +0.62%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.016.70"
$ws.Range("E43").Value = "  -1.63%  "

# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.24"
$ws.Range("E44").Value = "  +9.23%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "  +1.76%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'10.25"
$ws.Range("E46").Value = "  +2.04%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.91"
$ws.Range("E47").Value = "  +2.65%  "

# Row 48
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").Value = "'3.17"
$ws.Range("E48").Value = "  +216.47%  "

# Row 49
$ws.Range("E49").Value = "  +3.32%  "

# Row 50
$ws.Range("E50").Value = "  +0.58%  "

# Row 51
$ws.Range("D51").Value = "'72.48"
$ws.Range("E51").Value = "  +3.09%  "
